$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.518.98"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "'2.070.29"
$ws.Range("E3").Value = "  +5.61%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'236.01"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  +3.13%  "
$ws.Range("D7").Value = "'57.96"
$ws.Range("E7").Value = "  +9.17%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("D10").Value = "'57.69"
$ws.Range("D11").Value = "'0.0761"
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("E12").Value = "  +4.12%  "
$ws.Range("D13").Value = "'2.375.82"
$ws.Range("E13").Value = "  +5.54%  "
$ws.Range("D14").Value = "'14.27"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "'20.89"
$ws.Range("E15").Value = "  +5.21%  "
$ws.Range("D16").Value = "'0.777"
$ws.Range("E16").Value = "  +4.97%  "
$ws.Range("D17").Value = "'5.20"
$ws.Range("E17").Value = "  +4.53%  "
$ws.Range("D18").Value = "'2.075.46"
$ws.Range("E18").Value = "  +5.52%  "
$ws.Range("D19").Value = "'37.519.67"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").Value = "'6.16"
$ws.Range("E20").Value = "  +23.87%  "
$ws.Range("D21").Value = "'68.47"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "'0.0₃0811"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("D23").Value = "'224.49"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.44"
$ws.Range("E26").Value = "  +4.97%  "
$ws.Range("D27").Value = "'162.67"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").Value = "'8.83"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("E29").Value = "  +6.71%  "
$ws.Range("D30").Value = "'1.39"
$ws.Range("E30").Value = "  +7.63%  "
$ws.Range("D31").Value = "'19.29"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("D33").Value = "'0.0630"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("E34").Value = "  +15.81%  "
$ws.Range("D35").Value = "'4.46"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "'4.45"
$ws.Range("E36").Value = "  +6.96%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").Value = "'3.34"
$ws.Range("E39").Value = "  +5.04%  "
$ws.Range("D40").Value = "'5.85"
$ws.Range("E40").Value = "  +13.52%  "
$ws.Range("D41").Value = "'3.00"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").Value = "'0.0953"
$ws.Range("E42").Value = "  +9.60%  "
$ws.Range("D43").Value = "'1.476.01"
$ws.Range("E43").Value = "  +5.12%  "
$ws.Range("D44").Value = "'95.59"
$ws.Range("E44").Value = "  +10.65%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.31"
$ws.Range("E45").Value = "  +26.29%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0211"
$ws.Range("E46").Value = "  +5.25%  "
$ws.Range("D47").Value = "'16.11"
$ws.Range("E47").Value = "  +9.68%  "
$ws.Range("E48").Value = "  +4.31%  "
$ws.Range("D49").Value = "'7.31"
$ws.Range("E49").Value = "  +10.00%  "
$ws.Range("D50").Value = "'1.02"
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("D51").Value = "'2.92"
$ws.Range("E51").Value = "  +2.48%  "
